{"js": "// Revert \"Added support for importing devices from a csv file.\"\n//\n// Removes the two \"Import new devices\" bullet points (the heading bullet\n// and its csv-file sub-bullet) from the Management section, and carries\n// the \"_GoBack\" bookmark that used to sit at the end of that sub-bullet\n// over to the end of the following \"Notification\" heading, matching\n// upstream.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- locate the paragraphs we need by their text ---------------------\nlet importPara = null;\nlet csvPara = null;\nlet notificationPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"Import new devices\") {\n    importPara = para;\n  } else if (text.indexOf(\"Allow user to specify a csv file\") === 0) {\n    csvPara = para;\n  } else if (text === \"Notification\") {\n    notificationPara = para;\n  }\n}\n\n// --- remove both paragraphs (this also removes the old bookmark that\n//     lived at the end of the csv paragraph) --------------------------\nif (csvPara) {\n  csvPara.delete();\n}\nif (importPara) {\n  importPara.delete();\n}\nawait context.sync();\n\n// --- re-create the \"_GoBack\" bookmark at the end of the \"Notification\"\n//     heading paragraph (right after its text) ------------------------\nif (notificationPara) {\n  const endRange = notificationPara.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Revert \"Added support for importing devices from a csv file.\"\n#\n# Removes the two \"Import new devices\" bullet points (the heading bullet and\n# its csv-file sub-bullet) from the Management section, and carries the\n# \"_GoBack\" bookmark that used to sit at the end of that sub-bullet over to\n# the end of the following \"Notification\" heading, matching upstream.\n\n$d = $word.ActiveDocument\n\n# --- locate the two paragraphs to remove, by their text -------------------\n$importPara = $null\n$csvPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\")\n    if ($t -eq \"Import new devices\") {\n        $importPara = $p\n    } elseif ($t.StartsWith(\"Allow user to specify a csv file\")) {\n        $csvPara = $p\n    }\n}\n\n# --- remove both paragraphs (and whatever they contain, incl. the old\n#     bookmark) in a single range delete -----------------------------------\nif ($importPara -ne $null -and $csvPara -ne $null) {\n    $rng = $d.Range($importPara.Range.Start, $csvPara.Range.End)\n    $rng.Delete()\n}\n\n# --- re-create the \"_GoBack\" bookmark at the end of the \"Notification\"\n#     heading (right after the run text, before the paragraph mark) -------\n$notificationPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Notification\") {\n        $notificationPara = $p\n        break\n    }\n}\n\nif ($notificationPara -ne $null) {\n    $insPos = $notificationPara.Range.End - 1\n    $tmp = $d.Range($insPos, $insPos)\n    # Inserting a temporary marker char first avoids an edge-case where a\n    # zero-length range sitting exactly one char before a paragraph mark\n    # resolves to the wrong spot; bookmark the marker, then delete just the\n    # marker text, leaving the (now collapsed) bookmark behind in place.\n    $tmp.InsertBefore(\"x\")\n    $d.Bookmarks.Add(\"_GoBack\", $tmp)\n    $bm = $d.Bookmarks(\"_GoBack\")\n    $bm.Range.Text = \"\"\n}\n"}
